$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Overview ----
$ws = $wb.Worksheets.Item(1)

# Update cell values
$ws.Range('A1').Value = 'File Name'
$ws.Range('B1').Value = 'zh-cn'
$ws.Range('C1').Value = 'de-de'
$ws.Range('A2').Value = '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md'
$ws.Range('B2').Value = 'Handed back: in sync with en-US'
$ws.Range('C2').Value = 'Handed back: in sync with en-US'
$ws.Range('A3').Value = '340eefdd-c01c-4f44-96d8-19a1448a7aab.md'
$ws.Range('B3').Value = 'Handed back: in sync with en-US'
$ws.Range('C3').Value = 'Handed back: in sync with en-US'
$ws.Range('A4').Value = '.localization-config'
$ws.Range('B4').Value = 'Not to be localized'
$ws.Range('C4').Value = 'Not to be localized'

# Rebuild hyperlinks (same targets, refreshed display text) in original order
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/813e24a9bd2d64d7a165d7db9260070df0645306/e2e/340eefdd-c01c-4f44-96d8-19a1448a7aab.md', "", "", '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md', "", "", '340eefdd-c01c-4f44-96d8-19a1448a7aab.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/.localization-config', "", "", '.localization-config') | Out-Null

# ---- Sheet 2: zh-cn ----
$ws = $wb.Worksheets.Item(2)

# Update cell values
$ws.Range('A1').Value = 'Source File Name'
$ws.Range('B1').Value = 'Status'
$ws.Range('C1').Value = 'Latest Handoff File'
$ws.Range('D1').Value = 'Latest Handoff Datetime'
$ws.Range('E1').Value = 'Latest Target File'
$ws.Range('F1').Value = 'Latest Handback File'
$ws.Range('G1').Value = 'Latest Handback DateTime'
$ws.Range('H1').Value = 'Handoff Reason'
$ws.Range('I1').Value = 'Dependency From'
$ws.Range('A2').Value = '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md'
$ws.Range('B2').Value = 'Handed back: in sync with en-US'
$ws.Range('C2').Value = '2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf'
$ws.Range('D2').Value = '2016-01-25 03:38:37'
$ws.Range('E2').Value = '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md'
$ws.Range('F2').Value = '2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf'
$ws.Range('G2').Value = '2016-01-25 03:39:28'
$ws.Range('H2').Value = 'Include'
$ws.Range('A3').Value = '340eefdd-c01c-4f44-96d8-19a1448a7aab.md'
$ws.Range('B3').Value = 'Handed back: in sync with en-US'
$ws.Range('C3').Value = '340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.zh-cn.xlf'
$ws.Range('D3').Value = '2016-01-25 03:36:25'
$ws.Range('E3').Value = '340eefdd-c01c-4f44-96d8-19a1448a7aab.md'
$ws.Range('F3').Value = '340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.zh-cn.xlf'
$ws.Range('G3').Value = '2016-01-25 03:37:17'
$ws.Range('H3').Value = 'Include'
$ws.Range('A4').Value = '.localization-config'
$ws.Range('B4').Value = 'Not to be localized'
$ws.Range('D4').Value = '0001-01-01 00:00:00'
$ws.Range('G4').Value = '0001-01-01 00:00:00'
$ws.Range('H4').Value = 'Ignored'

# Rebuild hyperlinks (same targets, refreshed display text) in original order
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/813e24a9bd2d64d7a165d7db9260070df0645306/e2e/340eefdd-c01c-4f44-96d8-19a1448a7aab.md', "", "", '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0a0720ce4959c988dafb13560021c222c14a2f2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.zh-cn.xlf', "", "", '2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('E2'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ab1fbac9be8d4c9a6ac7d354d4d9b2b84a3b8a76/e2e/340eefdd-c01c-4f44-96d8-19a1448a7aab.md', "", "", '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/334cbdbd49308293623df4cbc2bc66558c405860/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/qimu/340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.zh-cn.xlf', "", "", '2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md', "", "", '340eefdd-c01c-4f44-96d8-19a1448a7aab.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/547e5119a31a5aa0c372863ce57cca36cea6165a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf', "", "", '340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('E3'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ab1fbac9be8d4c9a6ac7d354d4d9b2b84a3b8a76/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md', "", "", '340eefdd-c01c-4f44-96d8-19a1448a7aab.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/334cbdbd49308293623df4cbc2bc66558c405860/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf', "", "", '340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.zh-cn.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/.localization-config', "", "", '.localization-config') | Out-Null

# ---- Sheet 3: de-de ----
$ws = $wb.Worksheets.Item(3)

# Update cell values
$ws.Range('A1').Value = 'Source File Name'
$ws.Range('B1').Value = 'Status'
$ws.Range('C1').Value = 'Latest Handoff File'
$ws.Range('D1').Value = 'Latest Handoff Datetime'
$ws.Range('E1').Value = 'Latest Target File'
$ws.Range('F1').Value = 'Latest Handback File'
$ws.Range('G1').Value = 'Latest Handback DateTime'
$ws.Range('H1').Value = 'Handoff Reason'
$ws.Range('I1').Value = 'Dependency From'
$ws.Range('A2').Value = '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md'
$ws.Range('B2').Value = 'Handed back: in sync with en-US'
$ws.Range('C2').Value = '2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf'
$ws.Range('D2').Value = '2016-01-25 03:38:49'
$ws.Range('E2').Value = '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md'
$ws.Range('F2').Value = '2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf'
$ws.Range('G2').Value = '2016-01-25 03:39:44'
$ws.Range('H2').Value = 'Include'
$ws.Range('A3').Value = '340eefdd-c01c-4f44-96d8-19a1448a7aab.md'
$ws.Range('B3').Value = 'Handed back: in sync with en-US'
$ws.Range('C3').Value = '340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.de-de.xlf'
$ws.Range('D3').Value = '2016-01-25 03:36:37'
$ws.Range('E3').Value = '340eefdd-c01c-4f44-96d8-19a1448a7aab.md'
$ws.Range('F3').Value = '340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.de-de.xlf'
$ws.Range('G3').Value = '2016-01-25 03:37:33'
$ws.Range('H3').Value = 'Include'
$ws.Range('A4').Value = '.localization-config'
$ws.Range('B4').Value = 'Not to be localized'
$ws.Range('D4').Value = '0001-01-01 00:00:00'
$ws.Range('G4').Value = '0001-01-01 00:00:00'
$ws.Range('H4').Value = 'Ignored'

# Rebuild hyperlinks (same targets, refreshed display text) in original order
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/813e24a9bd2d64d7a165d7db9260070df0645306/e2e/340eefdd-c01c-4f44-96d8-19a1448a7aab.md', "", "", '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b6dc93023672594488cbc93afbe4ce0d2624122/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.de-de.xlf', "", "", '2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('E2'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c907cee7e7768369696dedf036c1157bcdf65e3f/e2e/340eefdd-c01c-4f44-96d8-19a1448a7aab.md', "", "", '2ac41cb1-a240-442c-a524-4efa10ef2ca6.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/c4ce757c93c4203ab47cd151dcec984755e5e479/ol-handback/OpenLocalizationTestOrg/oltest.de-de/qimu/340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.de-de.xlf', "", "", '2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md', "", "", '340eefdd-c01c-4f44-96d8-19a1448a7aab.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/414197c5ac1b0ab62b3841ee1b86f4ee1b7e736a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf', "", "", '340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('E3'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c907cee7e7768369696dedf036c1157bcdf65e3f/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md', "", "", '340eefdd-c01c-4f44-96d8-19a1448a7aab.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/c4ce757c93c4203ab47cd151dcec984755e5e479/ol-handback/OpenLocalizationTestOrg/oltest.de-de/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf', "", "", '340eefdd-c01c-4f44-96d8-19a1448a7aab.05a2c2c3a444464970b708d766c35f4ba88bb0ce.de-de.xlf') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/.localization-config', "", "", '.localization-config') | Out-Null

Write-Output 'done'
